$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7 (d=7)
$ws.Range("B7").Value = 96.97943286999613
$ws.Range("C7").Value = 97.8074541382358
$ws.Range("D7").Value = 98.54831033822251
$ws.Range("E7").Value = 98.0932548426664

# Row 8 (d=10)
$ws.Range("B8").Value = 97.459627063686
$ws.Range("C8").Value = 98.05797150468032
$ws.Range("D8").Value = 98.12388290307297
$ws.Range("E8").Value = 98.66986366616644

# Row 9 (last)
$ws.Range("B9").Value = 95.92516712784717
$ws.Range("C9").Value = 94.6367803025617
$ws.Range("D9").Value = 94.65867866342001
$ws.Range("E9").Value = 96.23645965195773
